# Add a "number_of_run" parameter column to the "scenarios" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at D (shifts existing D..H to E..I).
$ws.Columns.Item(4).Insert()

# Header + values for the new column.
$ws.Cells.Item(1, 4).Value = "number_of_run"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(4, 4).Value = 1

# Give the new header cell a distinct format (fill/border applied, no visible change).
$ws.Range("D1").Interior.ColorIndex = -4142

# Restore explicit column width for the new column D.
$ws.Columns.Item(4).ColumnWidth = 13.285714285714286

# Move the selection / active sheet to match the saved view state.
$ws.Select()
$ws.Range("D5").Select()

# pageSetup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
